$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "students" right before "partnership"
$beforeSheet = $wb.Worksheets.Item("partnership")
$students = $wb.Worksheets.Add($beforeSheet)
$students.Name = "students"

# Header row
$students.Range("A1").Value = "Year"
$students.Range("B1").Value = "Share"

# Data rows: year vs share
$data = @(
    @(2010, 0.239),
    @(2011, 0.239),
    @(2012, 0.234),
    @(2013, 0.23),
    @(2014, 0.226),
    @(2015, 0.219),
    @(2016, 0.218),
    @(2017, 0.213),
    @(2018, 0.21),
    @(2019, 0.205),
    @(2020, 0.204),
    @(2021, 0.197),
    @(2022, 0.195),
    @(2023, 0.188)
)

$row = 2
foreach ($pair in $data) {
    $students.Cells.Item($row, 1).Value = $pair[0]
    $students.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Make B2 the selected cell on the new sheet, and activate the sheet/tab
$students.Range("B2").Select() | Out-Null
$students.Activate() | Out-Null
